$wb = $excel.ActiveWorkbook

$wsProduct = $wb.Worksheets.Item("Productdata")
$wsProduct.Range("G2").Value = 70

# The H2:H11 cells are empty string-typed cells (no <v>) in the source file.
# Re-assert them as blank so the load/save round-trip does not resolve the
# missing value to shared-string index 0 ("Name"); this keeps those cells
# untouched/blank exactly as in the original workbook (the diff does not
# modify column H).
$wsProduct.Range("H2:H11").Value = ""

$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775
